# Add the missing Day-8 "Coding" journal row (row 58) on the Journal sheet,
# describing the implementation of Issue #10 (title screen / Move() rename
# work), and update the selected cell to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Row 58: Tache="Coding", Duree=50, Explications="Implementation of Issue #10
# to the program", Etats="In the work", Date de fin = 14:45
$ws.Range("A58").Value = "Coding"
$ws.Range("C58").Value = 50
$ws.Range("D58").Value = "Implementation of Issue #10 to the program"
$ws.Range("E58").Value = "In the work"
$ws.Range("F58").NumberFormat = "h:mm"
$ws.Range("F58").Value = 0.61458333333333337

# Move the active selection to match the saved view state.
$ws.Range("C59").Select()
